$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.360.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "'2.314.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.52%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'310.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.96%  "
$ws.Range("D6").Value = "'106.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.61%  "
$ws.Range("D7").Value = "'0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.40%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.75%  "
$ws.Range("D10").Value = "'40.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "'8.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.67%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("E14").Value = "  -3.57%  "
$ws.Range("D15").Value = "'15.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.78%  "
$ws.Range("D16").Value = "'2.664.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.49%  "
$ws.Range("D17").Value = "'2.329.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.67%  "
$ws.Range("D18").Value = "'42.296.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("E19").Value = "  -4.94%  "
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("E22").Value = "  -6.46%  "
$ws.Range("D23").Value = "'261.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.90%  "
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").Value = "'9.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.05%  "
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").Value = "'11.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.63%  "
$ws.Range("D28").Value = "'23.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.99%  "
$ws.Range("D30").Value = "'35.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").Value = "'165.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'0.0895"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("E35").Value = "  +12.17%  "
$ws.Range("E36").Value = "  -3.55%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("E38").Value = "  -2.08%  "
$ws.Range("D39").Value = "'3.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.91%  "
$ws.Range("D40").Value = "'2.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.63%  "
$ws.Range("D41").Value = "'100.45"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.63%  "
$ws.Range("E42").Value = "  -3.93%  "
$ws.Range("D43").Value = "'70.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("D44").Value = "'0.231"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'12.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.85%  "
$ws.Range("E47").Value = "  -5.02%  "
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "'9.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("D50").Value = "'73.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.83%  "
$ws.Range("D51").Value = "'1.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.55%  "
